$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill column E (duplicate_image_filename) with "NA" for rows 2 through 21
$ws.Range("E2:E21").Value = "NA"
